$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B3:B32 values (District heating demand per country)
$values = @(
    56.891999999999996,
    16.177,
    29.441000000000003,
    15.66,
    0.88500000000000001,
    103.03399999999999,
    478.25199999999995,
    106.63600000000001,
    20.350000000000001,
    3.1429999999999998,
    17.2,
    118.413,
    140.54899999999998,
    9.7460000000000004,
    29.109000000000002,
    0.61899999999999999,
    112.07000000000001,
    33.222000000000001,
    1.885,
    22.69,
    0.4,
    151.94999999999999,
    33.256,
    214.71400000000003,
    18.620999999999999,
    69.711000000000013,
    146.411,
    7.6950000000000003,
    26.805,
    45.171000000000006
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# B18 (IE) previously had no value and a plain (unformatted) style; now that it
# carries a number, give it the same numeric style used by the other data cells.
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(18, 2).Value = 0.61899999999999999
$excel.CutCopyMode = 0

# Update the selection to B3:B32 with active cell B3
$ws.Range("B3:B32").Select()
